$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value2 = "fullRNASeq"
    }
}
